$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '39.662.92'
$ws.Range("E2").Value = '  +1.32%  '

# Row 3
$ws.Range("D3").Value = '2.208.33'
$ws.Range("E3").Value = '  +0.61%  '

# Row 4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").Value = "'291.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.75%  '

# Row 6
$ws.Range("D6").Value = "'86.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.91%  '

# Row 7
$ws.Range("D7").Value = "'0.514"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.37%  '

# Row 8
$ws.Range("E8").Value = '  -0.10%  '

# Row 9
$ws.Range("E9").Value = '  +1.19%  '

# Row 10
$ws.Range("D10").Value = "'30.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.95%  '

# Row 11
$ws.Range("D11").Value = "'0.0785"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.50%  '

# Row 12
$ws.Range("D12").Value = "'47.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.63%  '

# Row 13
$ws.Range("E13").Value = '  +1.85%  '

# Row 14
$ws.Range("D14").Value = "'6.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.44%  '

# Row 15
$ws.Range("D15").Value = '2.549.22'
$ws.Range("E15").Value = '  +0.35%  '

# Row 16
$ws.Range("D16").Value = "'14.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.41%  '

# Row 17
$ws.Range("D17").Value = '2.217.06'
$ws.Range("E17").Value = '  +0.19%  '

# Row 18
$ws.Range("D18").Value = "'0.726"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.89%  '

# Row 19
$ws.Range("D19").Value = '39.628.21'
$ws.Range("E19").Value = '  +1.42%  '

# Row 20
$ws.Range("D20").Value = "'11.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +11.49%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0878'
$ws.Range("E21").Value = '  +1.44%  '

# Row 22
$ws.Range("E22").Value = '  +2.06%  '

# Row 23
$ws.Range("D23").Value = "'65.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.95%  '

# Row 24
$ws.Range("D24").Value = "'235.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.14%  '

# Row 25
$ws.Range("E25").Value = '  +0.01%  '

# Row 26
$ws.Range("E26").Value = '  +3.19%  '

# Row 27
$ws.Range("E27").Value = '  +2.55%  '

# Row 28
$ws.Range("D28").Value = "'22.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.32%  '

# Row 29
$ws.Range("E29").Value = '  +1.31%  '

# Row 30
$ws.Range("D30").Value = "'9.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.83%  '

# Row 31
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = "'32.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.10%  '

# Row 32
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = "'151.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.25%  '

# Row 33
$ws.Range("E33").Value = '  -0.13%  '

# Row 34
$ws.Range("D34").Value = "'4.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.42%  '

# Row 35
$ws.Range("D35").Value = "'0.0716"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.92%  '

# Row 36
$ws.Range("D36").Value = "'2.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.00%  '

# Row 37
$ws.Range("E37").Value = '  +2.58%  '

# Row 38
$ws.Range("E38").Value = '  +7.00%  '

# Row 39
$ws.Range("D39").Value = "'15.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.48%  '

# Row 40
$ws.Range("E40").Value = '  +3.86%  '

# Row 41
$ws.Range("E41").Value = '  +3.77%  '

# Row 42
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.071.18'
$ws.Range("E42").Value = '  +9.17%  '

# Row 43
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = "'3.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.84%  '

# Row 44
$ws.Range("E44").Value = '  +4.67%  '

# Row 45
$ws.Range("E45").Value = '  +3.85%  '

# Row 46
$ws.Range("D46").Value = "'9.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +10.75%  '

# Row 47
$ws.Range("D47").Value = "'17.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +12.38%  '

# Row 48
$ws.Range("E48").Value = '  +0.79%  '

# Row 49
$ws.Range("D49").Value = '2.418.67'
$ws.Range("E49").Value = '  +0.41%  '

# Row 50
$ws.Range("D50").Value = "'70.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.18%  '

# Row 51
$ws.Range("D51").Value = "'88.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.44%  '
